$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "Sede: {sede}" -> "Sede: {sede} - {city}"
#    The placeholder is stored as three bold runs: "{", "sede", "}".
#    Rewrite them back-to-front (right to left) so earlier offsets stay valid.
# ---------------------------------------------------------------------------
$anchor = $d.Content
$null = $anchor.Find.Execute("Sede: {sede}", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$base = $anchor.End - 6   # 6 = length of the literal text "{sede}"

$rThird = $d.Range($base + 5, $base + 6)     # "}"
$rThird.Text = "{city}"

$rSecond = $d.Range($base + 1, $base + 5)    # "sede"
$rSecond.Text = " - "

$rFirst = $d.Range($base, $base + 1)         # "{"
$rFirst.Text = "{sede}"

# ---------------------------------------------------------------------------
# 2) Collapse the "Se propone ... Unidades de Competencia:" sentence (which was
#    split across several runs) into one run with the same visible text, and
#    carry the paragraph's usual font onto the blank paragraph right after it.
# ---------------------------------------------------------------------------
$sentence = "Se propone la continuidad en la actual sede para ser asesorado/a en las siguientes Unidades de Competencia:"
$null = $d.Content.Find.Execute($sentence, $true, $false, $false, $false, $false, $true, 1, $false, $sentence, 2)

$paragraphs = $d.Paragraphs
for ($i = 1; $i -le $paragraphs.Count; $i++) {
    $para = $paragraphs.Item($i)
    if ($para.Range.Text.TrimEnd([char]13) -eq $sentence) {
        $blank = $paragraphs.Item($i + 1)
        $blank.Range.Font.NameBi = "Calibri"
        break
    }
}
